$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.769.32"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.652.96"
$ws.Range("E3").Value = "  +2.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.88"
$ws.Range("E5").Value = "  +0.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.04"
$ws.Range("E6").Value = "  +1.91%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.58"
$ws.Range("E9").Value = "  +1.92%  "

$ws.Range("E10").Value = "  +1.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.375"
$ws.Range("E11").Value = "  +2.82%  "

$ws.Range("E12").Value = "  +1.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.126.19"
$ws.Range("E13").Value = "  +2.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.69"
$ws.Range("E14").Value = "  +10.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.723.38"
$ws.Range("E15").Value = "  +0.41%  "

$ws.Range("E16").Value = "  +1.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.658.99"
$ws.Range("E17").Value = "  +1.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.55"
$ws.Range("E18").Value = "  +1.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.73"
$ws.Range("E19").Value = "  +1.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "350.56"
$ws.Range("E20").Value = "  +1.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.96"
$ws.Range("E21").Value = "  -0.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.534"
$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.02"
$ws.Range("E24").Value = "  +1.20%  "

$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.162"
$ws.Range("E26").Value = "  +1.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.17"
$ws.Range("E27").Value = "  +5.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.99"
$ws.Range("E28").Value = "  +9.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0810"
$ws.Range("E29").Value = "  +3.03%  "

$ws.Range("E30").Value = "  +5.46%  "

$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.18"
$ws.Range("E32").Value = "  +4.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.86"
$ws.Range("E33").Value = "  +1.98%  "

$ws.Range("E34").Value = "  +9.01%  "

$ws.Range("E35").Value = "  +5.24%  "

$ws.Range("E36").Value = "  +7.84%  "

$ws.Range("E37").Value = "  +2.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "327.36"
$ws.Range("E38").Value = "  +10.73%  "

$ws.Range("E39").Value = "  +4.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.36"
$ws.Range("E40").Value = "  +1.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.880"
$ws.Range("E41").Value = "  +3.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.25"
$ws.Range("E42").Value = "  +8.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.52"
$ws.Range("E43").Value = "  +4.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "135.08"
$ws.Range("E44").Value = "  -2.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0998"
$ws.Range("E45").Value = "  +1.41%  "

$ws.Range("E46").Value = "  +0.61%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.49"
$ws.Range("E47").Value = "  +3.38%  "

$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0558"
$ws.Range("E48").Value = "  +2.17%  "

$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.998"
$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0246"
$ws.Range("E50").Value = "  +2.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.134.21"
$ws.Range("E51").Value = "  +5.36%  "
